$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 971.7037
$ws.Range("I19").Value = 963.64703
$ws.Range("K19").Value = 963.64703
$ws.Range("M19").Value = -788.64703
$ws.Range("H43").Value = 29050.916
$ws.Range("I43").Value = 4624.0557
$ws.Range("J43").Value = 102331.5
$ws.Range("K43").Value = 4624.0557
$ws.Range("L43").Value = 102331.5
$ws.Range("M43").Value = -4555.0557
$ws.Range("N43").Value = -102469.5
$ws.Range("H58").Value = 1099.1666
$ws.Range("I58").Value = 1099.1666
$ws.Range("K58").Value = 3297.4998
$ws.Range("M58").Value = -3147.4998
$ws.Range("H103").Value = 1089.7
$ws.Range("I103").Value = 999.4
$ws.Range("J103").Value = 1180
$ws.Range("K103").Value = 2998.2
$ws.Range("L103").Value = 3540
$ws.Range("M103").Value = -2412.2
$ws.Range("N103").Value = -4712
$ws.Range("H116").Value = 5327.6665
$ws.Range("I116").Value = 4994.2
$ws.Range("K116").Value = 4994.2
$ws.Range("M116").Value = -1552.2
$ws.Range("H138").Value = 4031.37
$ws.Range("J138").Value = 4042.3264
$ws.Range("L138").Value = 12126.9792
$ws.Range("N138").Value = -22406.9792

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 23365.344
$ws.Range("I32").Value = 13472.292
$ws.Range("J32").Value = 59893.54
$ws.Range("K32").Value = 13472.292
$ws.Range("L32").Value = 59893.54
$ws.Range("M32").Value = -13185.292
$ws.Range("N32").Value = -60467.54
$ws.Range("H61").Value = 4656.722
$ws.Range("I61").Value = 3691.1904
$ws.Range("K61").Value = 3691.1904
$ws.Range("M61").Value = -3479.1904
$ws.Range("H74").Value = 1112.5454
$ws.Range("J74").Value = 369
$ws.Range("L74").Value = 369
$ws.Range("N74").Value = -2117
$ws.Range("H77").Value = 1112.5454
$ws.Range("J77").Value = 369
$ws.Range("L77").Value = 1845
$ws.Range("N77").Value = -10581
$ws.Range("H110").Value = 3049.7273
$ws.Range("I110").Value = 2854.7
$ws.Range("K110").Value = 2854.7
$ws.Range("M110").Value = -809.6999999999998
$ws.Range("H132").Value = 2809.9512
$ws.Range("I132").Value = 2447.2942
$ws.Range("J132").Value = 4571.4287
$ws.Range("K132").Value = 7341.882599999999
$ws.Range("L132").Value = 13714.2861
$ws.Range("M132").Value = -4811.882599999999
$ws.Range("N132").Value = -18774.2861
$ws.Range("H136").Value = 4656.722
$ws.Range("I136").Value = 3691.1904
$ws.Range("K136").Value = 11073.5712
$ws.Range("M136").Value = -8523.5712

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 1930.4667
$ws.Range("I20").Value = 1766
$ws.Range("J20").Value = 2999.5
$ws.Range("K20").Value = 1766
$ws.Range("L20").Value = 2999.5
$ws.Range("M20").Value = -1519
$ws.Range("N20").Value = -3493.5
$ws.Range("H99").Value = 3612.8
$ws.Range("I99").Value = 3612.8
$ws.Range("K99").Value = 3612.8
$ws.Range("M99").Value = -2114.8
$ws.Range("H105").Value = 4121
$ws.Range("J105").Value = 3828.9333
$ws.Range("L105").Value = 3828.9333
$ws.Range("N105").Value = -7322.933300000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 5190.9473
$ws.Range("I58").Value = 4808.6665
$ws.Range("K58").Value = 4808.6665
$ws.Range("M58").Value = -4605.6665
$ws.Range("H99").Value = 2768.7222
$ws.Range("I99").Value = 3604.625
$ws.Range("K99").Value = 3604.625
$ws.Range("M99").Value = -2106.625
$ws.Range("H107").Value = 6061040.5
$ws.Range("I107").Value = 6666914.5
$ws.Range("K107").Value = 6666914.5
$ws.Range("M107").Value = -6664994.5
$ws.Range("H126").Value = 2768.7222
$ws.Range("I126").Value = 3604.625
$ws.Range("K126").Value = 10813.875
$ws.Range("M126").Value = -8343.875
$ws.Range("H134").Value = 2206.3572
$ws.Range("I134").Value = 2089.2
$ws.Range("K134").Value = 6267.599999999999
$ws.Range("M134").Value = -3732.599999999999
$ws.Range("H136").Value = 5190.9473
$ws.Range("I136").Value = 4808.6665
$ws.Range("K136").Value = 14425.9995
$ws.Range("M136").Value = -11875.9995
$ws.Range("H138").Value = 53170.223
$ws.Range("J138").Value = 53170.223
$ws.Range("L138").Value = 53170.223
$ws.Range("N138").Value = -63450.223

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H132").Value = 2465.9429
$ws.Range("I132").Value = 1465.25
$ws.Range("J132").Value = 2762.4443
$ws.Range("K132").Value = 13187.25
$ws.Range("L132").Value = 24861.9987
$ws.Range("M132").Value = -10657.25
$ws.Range("N132").Value = -29921.9987
$ws.Range("H137").Value = 2714.6667
$ws.Range("I137").Value = 2403.4
$ws.Range("J137").Value = 2937
$ws.Range("K137").Value = 7210.200000000001
$ws.Range("L137").Value = 8811
$ws.Range("M137").Value = -2110.200000000001
$ws.Range("N137").Value = -19011

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 335.57144
$ws.Range("I2").Value = 262.83334
$ws.Range("K2").Value = 262.83334
$ws.Range("M2").Value = -149.83334
$ws.Range("H13").Value = 627.1429000000001
$ws.Range("J13").Value = 808
$ws.Range("L13").Value = 808
$ws.Range("N13").Value = -1086
$ws.Range("M40").ClearContents()
$ws.Range("H40").Value = 30018
$ws.Range("I40").Value = 0
$ws.Range("J40").Value = 30018
$ws.Range("K40").Value = 0
$ws.Range("L40").Value = 30018
$ws.Range("N40").Value = -30320
$ws.Range("H102").Value = 34201.32
$ws.Range("I102").Value = 37758
$ws.Range("K102").Value = 37758
$ws.Range("M102").Value = -36136
$ws.Range("N110").ClearContents()
$ws.Range("H110").Value = 45000
$ws.Range("J110").Value = 0
$ws.Range("L110").Value = 0
$ws.Range("H132").Value = 4569.3125
$ws.Range("I132").Value = 4807.3335
$ws.Range("K132").Value = 14422.0005
$ws.Range("M132").Value = -11892.0005

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 753.86664
$ws.Range("I55").Value = 747.53845
$ws.Range("J55").Value = 795
$ws.Range("K55").Value = 747.53845
$ws.Range("L55").Value = 795
$ws.Range("M55").Value = -574.53845
$ws.Range("N55").Value = -1141
$ws.Range("H61").Value = 2384.4285
$ws.Range("I61").Value = 2198.8333
$ws.Range("K61").Value = 2198.8333
$ws.Range("M61").Value = -1996.8333
$ws.Range("H113").Value = 2384.4285
$ws.Range("I113").Value = 2198.8333
$ws.Range("K113").Value = 2198.8333
$ws.Range("M113").Value = -28.83329999999978
$ws.Range("H130").Value = 79749.5
$ws.Range("J130").Value = 79749.5
$ws.Range("L130").Value = 79749.5
$ws.Range("N130").Value = -89789.5
$ws.Range("H132").Value = 3950.879
$ws.Range("J132").Value = 4051.4443
$ws.Range("L132").Value = 12154.3329
$ws.Range("N132").Value = -17214.3329

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H46").Value = 106176.86
$ws.Range("J46").Value = 106176.86
$ws.Range("L46").Value = 106176.86
$ws.Range("N46").Value = -106638.86
$ws.Range("H122").Value = 4925.5454
$ws.Range("I122").Value = 7698.3335
$ws.Range("K122").Value = 23095.0005
$ws.Range("M122").Value = -20645.0005
$ws.Range("H132").Value = 2978.25
$ws.Range("I132").Value = 2489.4285
$ws.Range("J132").Value = 6400
$ws.Range("K132").Value = 7468.2855
$ws.Range("L132").Value = 19200
$ws.Range("M132").Value = -4938.2855
$ws.Range("N132").Value = -24260
$ws.Range("H134").Value = 106176.86
$ws.Range("J134").Value = 106176.86
$ws.Range("L134").Value = 318530.58
$ws.Range("N134").Value = -323600.58
$ws.Range("H136").Value = 2577.38
$ws.Range("I136").Value = 2193.6897
$ws.Range("K136").Value = 6581.0691
$ws.Range("M136").Value = -4031.0691
$ws.Range("H138").Value = 72638.5
$ws.Range("I138").Value = 77777
$ws.Range("J138").Value = 67500
$ws.Range("K138").Value = 77777
$ws.Range("L138").Value = 67500
$ws.Range("M138").Value = -72637
$ws.Range("N138").Value = -77780
